$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert a new "room" column before the current column D.
# This pushes the existing "style" column (D) to E and the
# existing "link" column (E) to F, carrying their values/styles
# along (Excel's native column-insert behaviour).
# ------------------------------------------------------------------
$ws.Columns("D").Insert()

# ------------------------------------------------------------------
# Populate the new room column with header + values.
# ------------------------------------------------------------------
$ws.Range("D1").Value = "room"
$ws.Range("D2").Value = "phong-khach"
$ws.Range("D3").Value = "phong-ngu"
$ws.Range("D4").Value = "phong-bep"
$ws.Range("D5").Value = "phong-sinh-hoat-chung"
$ws.Range("D6").Value = "phong-khach"
$ws.Range("D7").Value = "phong-ngu"
$ws.Range("D8").Value = "phong-bep"

# ------------------------------------------------------------------
# The hyperlinks that used to live on column E now need to move to
# column F (the column-insert operation shifts cell content/styles
# but does not relocate the hyperlink metadata automatically).
# ------------------------------------------------------------------
$linkUrls = @(
    "https://sf-static.upanhlaylink.com/img/image_20251224be9afdd5a244e682829bd99b4236d340.jpg",
    "https://sf-static.upanhlaylink.com/img/image_20251224246168ce39e2a9578d9e4b23e32c8c4f.jpg",
    "https://sf-static.upanhlaylink.com/img/image_202512247639636fc96c515652b6bfd47abb897c.jpg",
    "https://sf-static.upanhlaylink.com/img/image_20251224ba2a52f4125656feae560b07ef1228c7.jpg",
    "https://sf-static.upanhlaylink.com/img/image_202512249e32b2e7884c7671185287f841718320.jpg",
    "https://sf-static.upanhlaylink.com/img/image_202512246a43fc6b96d548e6a8c3a6ef7fe5da55.jpg",
    "https://sf-static.upanhlaylink.com/img/image_20251224899b7011f28c4e9e8a36ab7e814e719f.jpg"
)

$ws.Range("E2:E8").Hyperlinks.Delete()

for ($i = 0; $i -lt $linkUrls.Count; $i++) {
    $row = $i + 2
    $target = $ws.Range("F$row")
    $ws.Hyperlinks.Add($target, $linkUrls[$i])
    # Re-assert the original Hyperlink cell style so the cell keeps
    # using the workbook's existing "Hyperlink" style slot instead of
    # a freshly minted duplicate.
    $target.Style = "Hyperlink"
}

# ------------------------------------------------------------------
# Minor column width tweak on column B.
# ------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 16.83

# ------------------------------------------------------------------
# Update the active selection shown in the saved worksheet view.
# ------------------------------------------------------------------
$ws.Range("D8").Select()
